$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data source corrected: columns J and K (rows 1-51) now hold a
# uniform pair of numeric constants instead of the old per-row values
# (J1/K1 used to be text "r"/"s"; J2:J51 was 0.5 and K2:K51 was 0.3).
$ws.Range("J1:J51").Value = 0.3
$ws.Range("K1:K51").Value = 0.6

# --- View state: selection moved to K1:K51 with K1 active, and the
# window was scrolled/zoomed to the new view (best-effort; not all
# window-state attributes round-trip through this host).
$excel.ActiveWindow.Zoom = 90
[void]$ws.Range("K1:K51").Select()
